# Updates cryptos list values (Price + Volume(1h)) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.583.17"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").Value = "2.316.53"
$ws.Range("E3").Value = "  -0.46%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("E8").Value = "  -0.64%  "

$ws.Range("D9").Value = "2.338.00"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.48%  "

$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.340"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.74%  "

$ws.Range("D15").Value = "2.733.87"
$ws.Range("E15").Value = "  -0.32%  "

$ws.Range("D16").Value = "56.583.33"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("D18").Value = "2.327.29"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "335.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("E28").Value = "  +1.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("D31").Value = "0.0₃0720"
$ws.Range("E31").Value = "  -2.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.43%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("E35").Value = "  -0.19%  "

$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.899"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.95%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "148.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.59%  "

$ws.Range("E42").Value = "  -1.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "286.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0928"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0500"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.559"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0215"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("E51").Value = "  -0.48%  "
